$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a string value to a cell while forcing text storage
# (avoids Excel auto-converting numeric-looking strings like "309.10"
# or "8.420" into floating point numbers, which would lose the exact
# textual representation). The NumberFormat is reset back to the
# worksheet's default ("Normal" style) afterwards so no new style index
# is left behind on the cell.
function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 41 <-> Row 42 content swap: FraxShare and TheSandbox switched rank
# positions, each also getting refreshed price / volume figures.
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws "D41" "0.5287"
$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D42" "7.154"
$ws.Range("E42").Value = "  -0.37%  "

# Refreshed Price / Volume(1h) values for the remaining rows
Set-TextValue $ws "D2" "26.913.89"
$ws.Range("E2").Value = "  +0.07%  "
Set-TextValue $ws "D3" "1.815.57"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue $ws "D5" "309.10"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  +0.08%  "
Set-TextValue $ws "D7" "0.4655"
$ws.Range("E7").Value = "  +1.13%  "
Set-TextValue $ws "D8" "0.3652"
$ws.Range("E8").Value = "  -1.36%  "
Set-TextValue $ws "D9" "0.07350"
$ws.Range("E9").Value = "  -0.34%  "
Set-TextValue $ws "D10" "0.8687"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("E11").Value = "  -1.21%  "
Set-TextValue $ws "D12" "1.819.10"
$ws.Range("E12").Value = "  -0.15%  "
Set-TextValue $ws "D13" "5.382"
$ws.Range("E13").Value = "  +0.46%  "
Set-TextValue $ws "D14" "0.07105"
$ws.Range("E14").Value = "  +0.92%  "
Set-TextValue $ws "D15" "6.505"
$ws.Range("E15").Value = "  -0.34%  "
Set-TextValue $ws "D16" "91.13"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("E17").Value = "  +0.19%  "
Set-TextValue $ws "D18" "0.000008680"
$ws.Range("E18").Value = "  -0.17%  "
Set-TextValue $ws "D20" "14.62"
$ws.Range("E20").Value = "  -0.63%  "
Set-TextValue $ws "D21" "26.937.60"
$ws.Range("E21").Value = "  +0.16%  "
Set-TextValue $ws "D22" "5.293"
$ws.Range("E22").Value = "  -0.67%  "
Set-TextValue $ws "D23" "10.56"
$ws.Range("E23").Value = "  -0.68%  "
Set-TextValue $ws "D24" "2.051.84"
$ws.Range("E24").Value = "  -0.57%  "
Set-TextValue $ws "D25" "1.894"
$ws.Range("E25").Value = "  -0.32%  "
Set-TextValue $ws "D26" "150.90"
$ws.Range("E26").Value = "  -0.26%  "
Set-TextValue $ws "D27" "18.32"
$ws.Range("E27").Value = "  -0.23%  "
Set-TextValue $ws "D28" "2.134"
$ws.Range("E28").Value = "  -0.30%  "
Set-TextValue $ws "D29" "5.253"
$ws.Range("E29").Value = "  -0.94%  "
Set-TextValue $ws "D30" "115.78"
$ws.Range("E30").Value = "  -0.14%  "
Set-TextValue $ws "D31" "0.08905"
$ws.Range("E31").Value = "  +0.08%  "
Set-TextValue $ws "D32" "0.7567"
$ws.Range("E32").Value = "  +0.59%  "
Set-TextValue $ws "D33" "1.163"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  +0.10%  "
Set-TextValue $ws "D37" "1.091"
$ws.Range("E37").Value = "  -0.76%  "
Set-TextValue $ws "D38" "0.05284"
$ws.Range("E38").Value = "  +0.95%  "
Set-TextValue $ws "D39" "0.01943"
$ws.Range("E39").Value = "  -1.43%  "
Set-TextValue $ws "D40" "2.965"
$ws.Range("E40").Value = "  +1.19%  "
Set-TextValue $ws "D43" "2.333"
$ws.Range("E43").Value = "  -3.88%  "
Set-TextValue $ws "D44" "0.1654"
$ws.Range("E44").Value = "  -0.61%  "
Set-TextValue $ws "D45" "8.420"
$ws.Range("E45").Value = "  -1.05%  "
Set-TextValue $ws "D46" "0.4847"
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("E48").Value = "  +0.08%  "
Set-TextValue $ws "D49" "103.15"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E50").Value = "  -0.83%  "
Set-TextValue $ws "D51" "0.06295"
$ws.Range("E51").Value = "  -0.01%  "
